# update DOCs for wk2 sprint2 (#22)
#
# Sheet1 ("BurnDown") tracks the sprint: column D is "Completed" items per
# day, column E ("Backlog") is the running remaining-work total that
# cascades off of D via shared formulas, and row 26 sums the D column.
#
# Sprint day 6 (row 7, 2019-04-19) had 1 item completed -> D7 = 1. Leave
# the selection on the cell that was just edited, matching the saved
# workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = 1

$ws.Range("D7").Select()
